# Applies the "Minor changes to presentation" edit:
#   - Slide 2 (Content Placeholder 2): swap the text of the 1st and 2nd
#     bullets ("Scrum Meetings" <-> "Changes from previous demo"), and
#     re-split all three of the first bullets' runs.
#   - Slide 3 title: "Changes from previous demo" -> "Changes from
#     previous demo and issues".
#   - Slide 3 (Content Placeholder 2): split "Technology Changes" into two
#     runs and add two new bullets, "Data sources" and "Internet".
#
# (The underlying chart's internal c:axId values also changed in the
# source diff, but those ids are an internal implementation detail that
# PowerPoint regenerates on its own when it resaves a chart part -- they
# are not exposed anywhere in the Chart/Axis COM object model, so there
# is nothing for a COM-interop script to set here.)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 2 - "overview": reshuffle the first three bullets
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(2)
$tr2 = $sh2.TextFrame.TextRange

# Bullet 1 was "Scrum Meetings" -> becomes "Changes from previous demo",
# typed/split as "Changes from previous " + "demo".
$para1 = $tr2.Paragraphs(1)
$para1.Text = "Changes from previous demo"
$split1 = "Changes from previous ".Length
$para1.Characters($split1 + 1, $para1.Length - $split1).Text = "demo"

# Bullet 2 was "Changes from previous demo" -> becomes "Scrum Meetings",
# split as "Scrum " + "Meetings".
$para2 = $tr2.Paragraphs(2)
$para2.Text = "Scrum Meetings"
$split2 = "Scrum ".Length
$para2.Characters($split2 + 1, $para2.Length - $split2).Text = "Meetings"

# Bullet 3 text stays "Client Meetings" but is re-split into "Client " +
# "Meetings".
$para3 = $tr2.Paragraphs(3)
$split3 = "Client ".Length
$para3.Characters($split3 + 1, $para3.Length - $split3).Text = "Meetings"

# ---------------------------------------------------------------------
# Slide 3 - title
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)

$title3 = $s3.Shapes.Item(1)
$trTitle = $title3.TextFrame.TextRange
$paraT = $trTitle.Paragraphs(1)
$paraT.Text = "Changes from previous demo and issues"
$splitT = "Changes from previous ".Length
$paraT.Characters($splitT + 1, $paraT.Length - $splitT).Text = "demo and issues"

# ---------------------------------------------------------------------
# Slide 3 - content placeholder: split "Technology Changes" and add the
# two new bullets "Data sources" / "Internet" right after it.
# ---------------------------------------------------------------------
$sh3 = $s3.Shapes.Item(2)
$tr3 = $sh3.TextFrame.TextRange

$paraTech = $tr3.Paragraphs(2)
$null = $paraTech.InsertAfter("`rData sources`rInternet")

# Re-fetch paragraph 2 (still "Technology Changes") now that the new
# paragraphs exist, and split it into two runs.
$paraTech2 = $tr3.Paragraphs(2)
$splitTech = "Technology ".Length
$paraTech2.Characters($splitTech + 1, $paraTech2.Length - $splitTech).Text = "Changes"
